# Regenerate save_data: update column G ("K") values for rows 2-35
# (commit: "regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 0
    3  = 4
    4  = 1
    5  = 0
    6  = 2
    7  = 2
    8  = 0
    9  = 2
    10 = 2
    11 = 1
    12 = 1
    13 = 1
    14 = 1
    15 = 0
    16 = 0
    17 = 0
    18 = 0
    19 = 1
    20 = 2
    21 = 1
    22 = 0
    23 = 0
    24 = 1
    25 = 2
    26 = 0
    27 = 1
    28 = 1
    29 = 2
    30 = 0
    31 = 1
    32 = 1
    33 = 2
    34 = 1
    35 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
